# Fixed network issues in DSL
# Update the "startStatusNotify with Synch" (VT293_0003) and "stopStatusNotify"
# (VT293_0005) test rows: replace the old screenshot-based steps/validation
# with the corrected wifi toggle + result validation steps, matching the
# pattern used by the other similar test rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")

# Row 4 - Objective: startStatusNotify with Synch (VT293_0003)
$ws.Range("G4").Value = "wait(3);`nvalidate1;`nlink_Click(network_test_link);`nvalidate2;`nSelectTestToRun(VT293_0003_string);`nClickRunTest(runtest_top_xpath);`nvalidate3;`nClickRunTest(runtest_bottom_xpath);`nwifi_Mode(OFF);`npress_Key(Back);`nvalidate3;`nvalidate4;`nwifi_Mode(ON);`nwait(2);`npress_Key(Back);`nvalidate3;`nvalidate5;"

$ws.Range("H4").Value = "validate1`n{`nvalidate_PageTitle=Manual specs`n};`nvalidate2`n{`nvalidate_PageTitle=Network JS Test`n};`nvalidate3`n{`nvalidate_Text_Exists=VT293-0003`n};`nvalidate4`n{`nvalidate_Result_notDisplayed=disconnected`n};`nvalidate5`n{`nvalidate_Result=Will fail`n};"

$ws.Rows.Item(4).RowHeight = 225.75

# Row 6 - Objective: stopStatusNotify (VT293_0005)
$ws.Range("G6").Value = "wait(3);`nvalidate1;`nlink_Click(network_test_link);`nvalidate2;`nSelectTestToRun(VT293_0005_string);`nClickRunTest(runtest_top_xpath);`nvalidate3;`nClickRunTest(runtest_bottom_xpath);`nwifi_Mode(OFF);`npress_Key(Back);`nvalidate3;`nvalidate4;`nwait(2);`nwifi_Mode(ON);`nwait(5);`npress_Key(Back);`nvalidate3;`nvalidate5;"

$ws.Range("H6").Value = "validate1`n{`nvalidate_PageTitle=Manual specs`n};`nvalidate2`n{`nvalidate_PageTitle=Network JS Test`n};`nvalidate3`n{`nvalidate_Text_Exists=VT293-0005`n};`nvalidate4`n{`nvalidate_Result=disconnected`n};`nvalidate5`n{`nvalidate_Result_notDisplayed=to connected`n};"
